# Apply the latest cryptocurrency price/volume-change refresh to the sheet.
# Column D ("Price") and column E ("Volume(1h)") are plain text cells in the
# original workbook (e.g. "30.702.69" is not a valid number), so for column D
# we force the cell to Text format before writing the new value and then
# restore the default "Normal" style, which keeps the value as text without
# leaving a stray number-format behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.701.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.74%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.121.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.67%  "

$ws.Range("E4").Value = "  +1.28%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.71%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5280"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.94%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4557"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.38%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.11"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.49%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09121"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.85%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.175"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.50"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.122.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.845"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.095"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.17%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.51%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001171"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.03%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.017"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06709"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.45%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.37%  "

$ws.Range("E21").Value = "  +1.27%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.446"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.776.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.83%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.98%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.383"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.61%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.369.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.92%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.70%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.32%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.555"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.34%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "135.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.58%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.211"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.61%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1079"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.51%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.420"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.638"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.18%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.959"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.53%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.57"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.991"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.88%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02669"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.67%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06885"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.14%  "

$ws.Range("E40").Value = "  +1.90%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.61%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6914"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.40%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.269"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.28"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.27%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6480"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.45%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.314"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.68%  "

$ws.Range("E47").Value = "  +16.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.708"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.74%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.260"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.90%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "83.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.40%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07320"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.73%  "
